$wb = $excel.ActiveWorkbook

# --- Sheet "A 1" (sheet1): add lat/lng header + data columns (F, G) ---
$ws1 = $wb.Worksheets.Item("A 1")
$ws1.Range("F1").Value = "lat"
$ws1.Range("G1").Value = "lng"
$ws1.Range("F2").Value = "lat"
$ws1.Range("G2").Value = "lng"

# --- Sheet "IdealHeaders" (sheet2): add lat/lng columns (H, I) ---
$ws2 = $wb.Worksheets.Item("IdealHeaders")
$ws2.Range("H1").Value = "lat"
$ws2.Range("I1").Value = "lng"

# Set the selection on IdealHeaders before switching away from it so the
# saved view state points at I2 (one row below the new last header).
[void]$ws2.Range("I2").Select()

# Make "A 1" the active sheet/tab and select G3, matching the saved view
# state captured in the edited workbook.
$ws1.Activate()
[void]$ws1.Range("G3").Select()
